$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.77588366666667
$ws.Range("H2").Value = 71.327651
$ws.Range("I2").Value = 0.201093431146956
$ws.Range("J2").Value = 0.2010934311469559
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 3996.959697900754
$ws.Range("R2").Value = 35972.63728110678
$ws.Range("S2").Value = 0.060009944943563
$ws.Range("T2").Value = 0.06000994494356301

$ws.Range("G3").Value = 23.77588366666667
$ws.Range("H3").Value = 71.327651
$ws.Range("I3").Value = 0.201093431146956
$ws.Range("J3").Value = 0.2010934311469559
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 3875.617296151918
$ws.Range("R3").Value = 34880.55566536726
$ws.Range("S3").Value = 0.058188122508853
$ws.Range("T3").Value = 0.05818812250885299

$ws.Range("G4").Value = 23.77588366666667
$ws.Range("H4").Value = 71.327651
$ws.Range("I4").Value = 0.201093431146956
$ws.Range("J4").Value = 0.2010934311469559
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 3946.642858699343
$ws.Range("R4").Value = 35519.78572829409
$ws.Range("S4").Value = 0.05925449305552005
$ws.Range("T4").Value = 0.05925449305552005

$ws.Range("G5").Value = 23.77588366666667
$ws.Range("H5").Value = 71.327651
$ws.Range("I5").Value = 0.201093431146956
$ws.Range("J5").Value = 0.2010934311469559
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 1574.599131132573
$ws.Range("R5").Value = 14171.39218019315
$ws.Range("S5").Value = 0.02364087063901992
$ws.Range("T5").Value = 0.02364087063901991

$ws.Range("I6").Value = 0.2703947904457373
$ws.Range("J6").Value = 0.2703947904457373
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 5374.402703110333
$ws.Range("R6").Value = 48369.624327993
$ws.Range("S6").Value = 0.08069073363125903
$ws.Range("T6").Value = 0.08069073363125905

$ws.Range("I7").Value = 0.2703947904457373
$ws.Range("J7").Value = 0.2703947904457373
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.07824106984734974
$ws.Range("T7").Value = 0.07824106984734974

$ws.Range("I8").Value = 0.2703947904457373
$ws.Range("J8").Value = 0.2703947904457373
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 5306.745539402114
$ws.Range("R8").Value = 47760.70985461902
$ws.Range("S8").Value = 0.07967493588095892
$ws.Range("T8").Value = 0.07967493588095892

$ws.Range("I9").Value = 0.2703947904457373
$ws.Range("J9").Value = 0.2703947904457373
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 2117.241720280215
$ws.Range("R9").Value = 19055.17548252194
$ws.Range("S9").Value = 0.03178805108616965
$ws.Range("T9").Value = 0.03178805108616965

$ws.Range("G10").Value = 14.51831366666667
$ws.Range("H10").Value = 43.554941
$ws.Range("I10").Value = 0.1227940694288843
$ws.Range("J10").Value = 0.1227940694288843
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 2440.671203674507
$ws.Range("R10").Value = 21966.04083307056
$ws.Range("S10").Value = 0.03664398833812899
$ws.Range("T10").Value = 0.036643988338129

$ws.Range("G11").Value = 14.51831366666667
$ws.Range("H11").Value = 43.554941
$ws.Range("I11").Value = 0.1227940694288843
$ws.Range("J11").Value = 0.1227940694288843
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 2366.575659031254
$ws.Range("R11").Value = 21299.18093128129
$ws.Range("S11").Value = 0.03553152539362139
$ws.Range("T11").Value = 0.03553152539362139

$ws.Range("G12").Value = 14.51831366666667
$ws.Range("H12").Value = 43.554941
$ws.Range("I12").Value = 0.1227940694288843
$ws.Range("J12").Value = 0.1227940694288843
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 2409.946135177243
$ws.Range("R12").Value = 21689.51521659519
$ws.Range("S12").Value = 0.03618268529574997
$ws.Range("T12").Value = 0.03618268529574997

$ws.Range("G13").Value = 14.51831366666667
$ws.Range("H13").Value = 43.554941
$ws.Range("I13").Value = 0.1227940694288843
$ws.Range("J13").Value = 0.1227940694288843
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 961.5005021703359
$ws.Range("R13").Value = 8653.504519533024
$ws.Range("S13").Value = 0.01443587040138396
$ws.Range("T13").Value = 0.01443587040138396

$ws.Range("G14").Value = 47.96923
$ws.Range("H14").Value = 143.90769
$ws.Range("I14").Value = 0.4057177089784224
$ws.Range("J14").Value = 0.4057177089784224
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 8064.098972612954
$ws.Range("R14").Value = 72576.89075351658
$ws.Range("S14").Value = 0.1210735588903009
$ws.Range("T14").Value = 0.1210735588903009

$ws.Range("G15").Value = 47.96923
$ws.Range("H15").Value = 143.90769
$ws.Range("I15").Value = 0.4057177089784224
$ws.Range("J15").Value = 0.4057177089784224
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 7819.283610128537
$ws.Range("R15").Value = 70373.55249115684
$ws.Range("S15").Value = 0.1173979260257153
$ws.Range("T15").Value = 0.1173979260257153

$ws.Range("G16").Value = 47.96923
$ws.Range("H16").Value = 143.90769
$ws.Range("I16").Value = 0.4057177089784224
$ws.Range("J16").Value = 0.4057177089784224
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 7962.5818190819
$ws.Range("R16").Value = 71663.2363717371
$ws.Range("S16").Value = 0.1195493907088141
$ws.Range("T16").Value = 0.1195493907088141

$ws.Range("G17").Value = 47.96923
$ws.Range("H17").Value = 143.90769
$ws.Range("I17").Value = 0.4057177089784224
$ws.Range("J17").Value = 0.4057177089784224
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 3176.845451384564
$ws.Range("R17").Value = 28591.60906246107
$ws.Range("S17").Value = 0.04769683335359214
$ws.Range("T17").Value = 0.04769683335359214
